# Edit: fill in results for "3ª Giornata lega / 6ª Giornata serie a" matches
# (rows 10-13) and rename team "Atletico Manontroppo" -> "SamPDORRR FC"
# everywhere it appears on the schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the team throughout the whole sheet (it appears many times across
# the season calendar) - mirrors a shared-string rename.
[void]$ws.Cells.Replace("Atletico Manontroppo", "SamPDORRR FC")

# Row 10: asdMalerba vs SamPDORRR FC -> 66 - 73.5, result 1-2
$ws.Range("B10").Value = 66
$ws.Range("C10").Value = 73.5
$ws.Range("E10").Value = "1-2"

# Row 11: Civediamoamaggio vs Capitan Martella -> 78.5 - 62.5, result 4-0
$ws.Range("B11").Value = 78.5
$ws.Range("C11").Value = 62.5
$ws.Range("E11").Value = "4-0"

# Row 12: Barriera-team vs Gargantua -> 70.5 - 77.5, result 2-3
$ws.Range("B12").Value = 70.5
$ws.Range("C12").Value = 77.5
$ws.Range("E12").Value = "2-3"

# Row 13: Paris FC vs CACCOLA F.C. -> 70.5 - 77.5, result 2-3
$ws.Range("B13").Value = 70.5
$ws.Range("C13").Value = 77.5
$ws.Range("E13").Value = "2-3"
